# Refresh the cryptos price/volume table (GitHub Actions data pull).
# Numeric-looking price strings are entered with a leading apostrophe so
# Excel keeps them as text (matching the sheet's original inline-string
# cells) instead of silently converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "98.757.26"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").Value = "3.311.32"
$ws.Range("E3").Value = "  -1.49%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'255.05"
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("D6").Value = "'629.88"
$ws.Range("E6").Value = "  +1.14%  "
$ws.Range("E7").Value = "  +21.78%  "
$ws.Range("D8").Value = "'0.412"
$ws.Range("E8").Value = "  +6.69%  "
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("D10").Value = "'1.01"
$ws.Range("E10").Value = "  +24.69%  "
$ws.Range("D11").Value = "3.308.11"
$ws.Range("E11").Value = "  -1.44%  "
$ws.Range("E12").Value = "  +3.23%  "
$ws.Range("D13").Value = "'42.46"
$ws.Range("E13").Value = "  +19.17%  "
$ws.Range("D14").Value = "98.426.70"
$ws.Range("E14").Value = "  +0.19%  "
$ws.Range("E15").Value = "  +1.94%  "
$ws.Range("D16").Value = "3.942.12"
$ws.Range("E16").Value = "  -1.35%  "
$ws.Range("E17").Value = "  -1.68%  "
$ws.Range("D18").Value = "3.312.42"
$ws.Range("E18").Value = "  -1.73%  "
$ws.Range("D19").Value = "'15.92"
$ws.Range("E19").Value = "  +7.31%  "
$ws.Range("D20").Value = "'3.49"
$ws.Range("E20").Value = "  -4.64%  "
$ws.Range("D21").Value = "'6.49"
$ws.Range("E21").Value = "  +9.56%  "
$ws.Range("D22").Value = "'489.39"
$ws.Range("E22").Value = "  +0.70%  "
$ws.Range("E23").Value = "  +2.55%  "
$ws.Range("E24").Value = "  -3.19%  "
$ws.Range("D25").Value = "'5.84"
$ws.Range("E25").Value = "  +2.40%  "
$ws.Range("D26").Value = "'0.348"
$ws.Range("E26").Value = "  +37.10%  "
$ws.Range("E27").Value = "  +1.45%  "
$ws.Range("D28").Value = "'12.24"
$ws.Range("E28").Value = "  +2.13%  "
$ws.Range("D29").Value = "3.489.46"
$ws.Range("E29").Value = "  -1.54%  "
$ws.Range("D30").Value = "'0.150"
$ws.Range("E30").Value = "  +19.67%  "
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("D32").Value = "'0.192"
$ws.Range("E32").Value = "  +3.16%  "
$ws.Range("D33").Value = "'10.95"
$ws.Range("E33").Value = "  +18.84%  "
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").Value = "'28.27"
$ws.Range("E35").Value = "  +2.99%  "
$ws.Range("D36").Value = "'0.483"
$ws.Range("E36").Value = "  +7.72%  "
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").Value = "'0.151"
$ws.Range("E37").Value = "  -1.02%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D38").Value = "'7.36"
$ws.Range("E38").Value = "  -0.03%  "
$ws.Range("E39").Value = "  +1.01%  "
$ws.Range("D40").Value = "'496.16"
$ws.Range("E40").Value = "  -5.06%  "
$ws.Range("B41").Value = "WhiteBITCoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D41").Value = "'24.73"
$ws.Range("E41").Value = "  -0.30%  "
$ws.Range("B42").Value = "MantraDAO"
$ws.Range("C42").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D42").Value = "'3.88"
$ws.Range("E42").Value = "  +2.31%  "
$ws.Range("E43").Value = "  -1.58%  "
$ws.Range("E44").Value = "  +1.45%  "
$ws.Range("E46").Value = "  -2.73%  "
$ws.Range("D47").Value = "'160.46"
$ws.Range("E47").Value = "  -0.47%  "
$ws.Range("D48").Value = "'1.97"
$ws.Range("E48").Value = "  +2.22%  "
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").Value = "'7.41"
$ws.Range("E49").Value = "  +15.70%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").Value = "'0.859"
$ws.Range("E50").Value = "  +8.08%  "
$ws.Range("D51").Value = "'4.80"
$ws.Range("E51").Value = "  +5.68%  "
